$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.626.83"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.891.53"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.54"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4905"
$ws.Range("E7").Value = "  +1.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2940"
$ws.Range("E8").Value = "  +1.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06707"
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.897.43"
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.03"
$ws.Range("E11").Value = "  +1.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07346"
$ws.Range("E12").Value = "  +2.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.162"
$ws.Range("E13").Value = "  +3.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.99"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6682"
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.581.01"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007871"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("E18").Value = "  +3.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.144.50"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.362"
$ws.Range("E21").Value = "  +13.05%  "
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "189.36"
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.202"
$ws.Range("E24").Value = "  +3.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.532"
$ws.Range("E25").Value = "  +3.44%  "
$ws.Range("E26").Value = "  +4.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.48"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.931"
$ws.Range("E28").Value = "  +4.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.465"
$ws.Range("E29").Value = "  +3.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.409"
$ws.Range("E30").Value = "  +4.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09158"
$ws.Range("E31").Value = "  +2.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.044"
$ws.Range("E32").Value = "  +3.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05249"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("E34").Value = "  +2.32%  "
$ws.Range("E35").Value = "  +2.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.727"
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01824"
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.696"
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9153"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.065"
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "75.36"
$ws.Range("E41").Value = "  +32.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4420"
$ws.Range("E42").Value = "  +2.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.931"
$ws.Range("E43").Value = "  +6.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "106.11"
$ws.Range("E44").Value = "  +2.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9931"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("E46").Value = "  +3.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.557"
$ws.Range("E47").Value = "  +3.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.50"
$ws.Range("E48").Value = "  +6.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.025"
$ws.Range("E49").Value = "  +3.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05842"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("E51").Value = "  +2.12%  "
